{"js": "// Replace each math-expression cell's text in the document's (single) table,\n// in row-major order, per the old->new mapping derived from the target diff.\n// Every \"old\" value is unique in the document and none of the \"new\" values\n// collide with any \"old\" value, so a direct positional assignment is safe.\nconst replacements = [\n  [\"6+75=\", \"4+77=\"],\n  [\"44-12=\", \"67+22=\"],\n  [\"31+11=\", \"75+1=\"],\n  [\"64-49=\", \"21+68=\"],\n  [\"20+44=\", \"55-10=\"],\n  [\"41+42=\", \"14+23=\"],\n  [\"94-84=\", \"68-8=\"],\n  [\"71+8=\", \"93-83=\"],\n  [\"30-12=\", \"7+89=\"],\n  [\"51-42=\", \"27+6=\"],\n  [\"34-2=\", \"7+21=\"],\n  [\"29+12=\", \"27-13=\"],\n  [\"29+13=\", \"21+74=\"],\n  [\"98-38=\", \"52+2=\"],\n  [\"1+20=\", \"39+10=\"],\n  [\"73+1=\", \"53-29=\"],\n  [\"70-65=\", \"70+17=\"],\n  [\"86-51=\", \"22+40=\"],\n  [\"61-45=\", \"98-30=\"],\n  [\"34+39=\", \"18+0=\"],\n  [\"43+56=\", \"47-41=\"],\n  [\"74-56=\", \"35+60=\"],\n  [\"32-7=\", \"70+24=\"],\n  [\"9+72=\", \"46+17=\"],\n  [\"26+68=\", \"71-1=\"],\n  [\"57-15=\", \"47+21=\"],\n  [\"69-3=\", \"52-17=\"],\n  [\"82-0=\", \"87-49=\"],\n  [\"4+54=\", \"48-47=\"],\n  [\"93+5=\", \"73-47=\"],\n  [\"6-0=\", \"11+61=\"],\n  [\"56-18=\", \"53-43=\"],\n  [\"65-41=\", \"55-39=\"],\n  [\"39-4=\", \"88+11=\"],\n  [\"85-9=\", \"2+87=\"],\n  [\"20+73=\", \"13-12=\"],\n  [\"75+15=\", \"87-75=\"],\n  [\"45-34=\", \"80-39=\"],\n  [\"3+86=\", \"71-13=\"],\n  [\"6-3=\", \"15+53=\"],\n  [\"47+7=\", \"23+33=\"],\n  [\"67-22=\", \"86-30=\"],\n  [\"61-5=\", \"66+2=\"],\n  [\"69+24=\", \"0+6=\"],\n  [\"94-38=\", \"60-33=\"],\n  [\"51+14=\", \"53+30=\"],\n  [\"57+21=\", \"45+47=\"],\n  [\"37+32=\", \"46-14=\"],\n  [\"51-30=\", \"90-24=\"],\n  [\"25+27=\", \"76-14=\"],\n  [\"81+14=\", \"29+11=\"],\n  [\"23+53=\", \"70-8=\"],\n  [\"67-14=\", \"3+61=\"],\n  [\"50+14=\", \"58+4=\"],\n  [\"42-33=\", \"66+33=\"],\n  [\"40-39=\", \"73+7=\"],\n  [\"70-66=\", \"47-39=\"],\n  [\"8+58=\", \"58-29=\"],\n  [\"94-50=\", \"79-9=\"],\n  [\"7+46=\", \"21-17=\"],\n  [\"55-31=\", \"41-38=\"],\n  [\"50+15=\", \"76-41=\"],\n  [\"75-1=\", \"22+23=\"],\n  [\"61-42=\", \"69-62=\"],\n  [\"5+42=\", \"42-15=\"],\n  [\"33+4=\", \"55-1=\"],\n  [\"82-9=\", \"18+55=\"],\n  [\"76-9=\", \"65-11=\"],\n  [\"40+22=\", \"9+87=\"],\n  [\"30+41=\", \"74+12=\"],\n  [\"99-25=\", \"59-39=\"],\n  [\"6+88=\", \"89-14=\"],\n  [\"8+48=\", \"31+43=\"],\n  [\"35+11=\", \"69+1=\"],\n  [\"14+77=\", \"2+81=\"],\n  [\"56+18=\", \"95-50=\"],\n  [\"85-1=\", \"10+5=\"],\n  [\"99-75=\", \"6+24=\"],\n  [\"52-52=\", \"51+30=\"],\n  [\"17-15=\", \"90-12=\"],\n  [\"87-76=\", \"51-24=\"],\n  [\"99-95=\", \"63+10=\"],\n  [\"75-22=\", \"37+61=\"],\n  [\"26+51=\", \"20+52=\"],\n  [\"27-9=\", \"67-21=\"],\n  [\"78-71=\", \"67+14=\"],\n  [\"32+10=\", \"43+42=\"],\n  [\"11+1=\", \"16-16=\"],\n  [\"94-29=\", \"30+5=\"],\n  [\"44-30=\", \"37+38=\"],\n  [\"97-8=\", \"92-13=\"],\n  [\"45-45=\", \"98-29=\"],\n  [\"5+87=\", \"33+28=\"],\n  [\"45+54=\", \"2+61=\"],\n  [\"92-75=\", \"31-29=\"],\n  [\"48-41=\", \"11+79=\"],\n  [\"33+35=\", \"39+60=\"],\n  [\"64-13=\", \"83-77=\"],\n  [\"99-39=\", \"94-34=\"],\n  [\"45-30=\", \"64-52=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n}\nawait context.sync();\n\n// Flatten every cell proxy in row-major (document) order.\nconst allCells = [];\nfor (const row of rows.items) {\n  for (const cell of row.cells.items) {\n    allCells.push(cell);\n  }\n}\n\nfor (const cell of allCells) {\n  cell.load(\"value\");\n}\nawait context.sync();\n\nif (allCells.length !== replacements.length) {\n  throw new Error(\n    `Expected ${replacements.length} table cells, found ${allCells.length}`\n  );\n}\n\nfor (let i = 0; i < allCells.length; i++) {\n  const [oldText, newText] = replacements[i];\n  const current = allCells[i].value;\n  if (current !== oldText) {\n    throw new Error(\n      `Cell ${i}: expected \"${oldText}\" but found \"${current}\"`\n    );\n  }\n  allCells[i].value = newText;\n}\n\nawait context.sync();\n", "ps1": "# Replace each math-expression cell's text in the document's (single) table,\n# in row-major order, per the old->new mapping derived from the target diff.\n# Every \"old\" value is unique in the document and none of the \"new\" values\n# collide with any \"old\" value, so a direct positional assignment is safe.\n$replacements = @(\n  @(\"6+75=\", \"4+77=\"),\n  @(\"44-12=\", \"67+22=\"),\n  @(\"31+11=\", \"75+1=\"),\n  @(\"64-49=\", \"21+68=\"),\n  @(\"20+44=\", \"55-10=\"),\n  @(\"41+42=\", \"14+23=\"),\n  @(\"94-84=\", \"68-8=\"),\n  @(\"71+8=\", \"93-83=\"),\n  @(\"30-12=\", \"7+89=\"),\n  @(\"51-42=\", \"27+6=\"),\n  @(\"34-2=\", \"7+21=\"),\n  @(\"29+12=\", \"27-13=\"),\n  @(\"29+13=\", \"21+74=\"),\n  @(\"98-38=\", \"52+2=\"),\n  @(\"1+20=\", \"39+10=\"),\n  @(\"73+1=\", \"53-29=\"),\n  @(\"70-65=\", \"70+17=\"),\n  @(\"86-51=\", \"22+40=\"),\n  @(\"61-45=\", \"98-30=\"),\n  @(\"34+39=\", \"18+0=\"),\n  @(\"43+56=\", \"47-41=\"),\n  @(\"74-56=\", \"35+60=\"),\n  @(\"32-7=\", \"70+24=\"),\n  @(\"9+72=\", \"46+17=\"),\n  @(\"26+68=\", \"71-1=\"),\n  @(\"57-15=\", \"47+21=\"),\n  @(\"69-3=\", \"52-17=\"),\n  @(\"82-0=\", \"87-49=\"),\n  @(\"4+54=\", \"48-47=\"),\n  @(\"93+5=\", \"73-47=\"),\n  @(\"6-0=\", \"11+61=\"),\n  @(\"56-18=\", \"53-43=\"),\n  @(\"65-41=\", \"55-39=\"),\n  @(\"39-4=\", \"88+11=\"),\n  @(\"85-9=\", \"2+87=\"),\n  @(\"20+73=\", \"13-12=\"),\n  @(\"75+15=\", \"87-75=\"),\n  @(\"45-34=\", \"80-39=\"),\n  @(\"3+86=\", \"71-13=\"),\n  @(\"6-3=\", \"15+53=\"),\n  @(\"47+7=\", \"23+33=\"),\n  @(\"67-22=\", \"86-30=\"),\n  @(\"61-5=\", \"66+2=\"),\n  @(\"69+24=\", \"0+6=\"),\n  @(\"94-38=\", \"60-33=\"),\n  @(\"51+14=\", \"53+30=\"),\n  @(\"57+21=\", \"45+47=\"),\n  @(\"37+32=\", \"46-14=\"),\n  @(\"51-30=\", \"90-24=\"),\n  @(\"25+27=\", \"76-14=\"),\n  @(\"81+14=\", \"29+11=\"),\n  @(\"23+53=\", \"70-8=\"),\n  @(\"67-14=\", \"3+61=\"),\n  @(\"50+14=\", \"58+4=\"),\n  @(\"42-33=\", \"66+33=\"),\n  @(\"40-39=\", \"73+7=\"),\n  @(\"70-66=\", \"47-39=\"),\n  @(\"8+58=\", \"58-29=\"),\n  @(\"94-50=\", \"79-9=\"),\n  @(\"7+46=\", \"21-17=\"),\n  @(\"55-31=\", \"41-38=\"),\n  @(\"50+15=\", \"76-41=\"),\n  @(\"75-1=\", \"22+23=\"),\n  @(\"61-42=\", \"69-62=\"),\n  @(\"5+42=\", \"42-15=\"),\n  @(\"33+4=\", \"55-1=\"),\n  @(\"82-9=\", \"18+55=\"),\n  @(\"76-9=\", \"65-11=\"),\n  @(\"40+22=\", \"9+87=\"),\n  @(\"30+41=\", \"74+12=\"),\n  @(\"99-25=\", \"59-39=\"),\n  @(\"6+88=\", \"89-14=\"),\n  @(\"8+48=\", \"31+43=\"),\n  @(\"35+11=\", \"69+1=\"),\n  @(\"14+77=\", \"2+81=\"),\n  @(\"56+18=\", \"95-50=\"),\n  @(\"85-1=\", \"10+5=\"),\n  @(\"99-75=\", \"6+24=\"),\n  @(\"52-52=\", \"51+30=\"),\n  @(\"17-15=\", \"90-12=\"),\n  @(\"87-76=\", \"51-24=\"),\n  @(\"99-95=\", \"63+10=\"),\n  @(\"75-22=\", \"37+61=\"),\n  @(\"26+51=\", \"20+52=\"),\n  @(\"27-9=\", \"67-21=\"),\n  @(\"78-71=\", \"67+14=\"),\n  @(\"32+10=\", \"43+42=\"),\n  @(\"11+1=\", \"16-16=\"),\n  @(\"94-29=\", \"30+5=\"),\n  @(\"44-30=\", \"37+38=\"),\n  @(\"97-8=\", \"92-13=\"),\n  @(\"45-45=\", \"98-29=\"),\n  @(\"5+87=\", \"33+28=\"),\n  @(\"45+54=\", \"2+61=\"),\n  @(\"92-75=\", \"31-29=\"),\n  @(\"48-41=\", \"11+79=\"),\n  @(\"33+35=\", \"39+60=\"),\n  @(\"64-13=\", \"83-77=\"),\n  @(\"99-39=\", \"94-34=\"),\n  @(\"45-30=\", \"64-52=\"),\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\nif (($rows * $cols) -ne $replacements.Count) {\n  throw \"Expected $($replacements.Count) table cells, found $($rows * $cols)\"\n}\n\n$k = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $pair = $replacements[$k]\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $cell = $t.Cell($r, $c)\n    $current = $cell.Range.Text\n    # Range.Text carries the trailing cell-mark pair (\"`r`a\"); strip it\n    # before comparing against the plain expected string.\n    $current = $current.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n      throw \"Cell ($r,$c): expected [$oldText] but found [$current]\"\n    }\n    $cell.Range.Text = $newText\n    $k = $k + 1\n  }\n}\n"}
